$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 7
    3 = 2
    4 = 10
    5 = 4
    6 = 3
    7 = 6
    8 = 7
    9 = 9
    10 = 13
    11 = 9
    12 = 11
    13 = 8
    14 = 7
    15 = 8
    16 = 10
    17 = 6
    18 = 8
    19 = 9
    20 = 8
    21 = 7
    22 = 8
    23 = 0
    24 = 9
    25 = 10
    26 = 9
    27 = 8
    28 = 7
    29 = 18
    30 = 9
    31 = 5
    32 = 8
    33 = 10
    34 = 5
    35 = 9
    36 = 6
    37 = 7
    38 = 5
    39 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
